# Blog workbook update: log a new "Troubleshooting" activity on 2023-07-16
# (row 20) for Paolo & Woldy, resolving a DB connection issue reported by
# Olivier, and move the sheet selection to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: new activity entry -------------------------------------------
$ws.Range("A20").Value = "7/16/2023"   # Date
$ws.Range("B20").Value = 1             # Paolo
$ws.Range("C20").Value = 0             # Yevhen
$ws.Range("D20").Value = 1             # Woldy
$ws.Range("E20").Value = 0             # Oliver
$ws.Range("F20").Value = 0             # Ke
$ws.Range("G20").Value = "Troubleshooting"   # type of activity
$ws.Range("H20").Value = 2                    # No. of hours spent

# Output ("resolved") is written before Purpose so the new shared-string
# table entries land in the same order as the authored workbook.
$ws.Range("J20").Value = "resolved"
$ws.Range("I20").Value = "DB connection issues with front-end for Olivier"

# Column I (Purpose) needs to widen to fit the new long entry - mirror the
# width Excel settled on after auto-fitting the column to its content.
$ws.Columns("I").ColumnWidth = 42.6

# --- Selection / view -------------------------------------------------------
$ws.Range("I15").Select()
